# Apply updated "results with all iterations" values to both sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet: "Full results" ---
$ws1 = $wb.Worksheets.Item("Full results")

$ws1.Cells.Item(2, 8).Value  = 0.588632892052824    # H2 completeind
$ws1.Cells.Item(2, 9).Value  = 0.28042484297219     # I2 completefam
$ws1.Cells.Item(2, 15).Value = 0.411420401256562    # O2 IORAD

$ws1.Cells.Item(3, 6).Value  = 0.603695825461784    # F3 condind
$ws1.Cells.Item(3, 7).Value  = 0.312032917347289    # G3 condfam

$ws1.Cells.Item(4, 3).Value  = 0.637314387588559    # C4 emptyind
$ws1.Cells.Item(4, 4).Value  = 0.362776158086459    # D4 emptyfam
$ws1.Cells.Item(4, 5).Value  = 1.00009054567502     # E4 totalvar
$ws1.Cells.Item(4, 10).Value = 0.362743313229904    # J4 Sibcorr
$ws1.Cells.Item(4, 11).Value = 0.312004666656564    # K4 condcorr
$ws1.Cells.Item(4, 12).Value = 0.0150615696308235   # L4 w
$ws1.Cells.Item(4, 13).Value = 0.0486770880266578   # M4 v
$ws1.Cells.Item(4, 14).Value = 0.327066236287387    # N4 IOLIB

# --- Sheet: "For plotting" ---
$ws2 = $wb.Worksheets.Item("For plotting")

$ws2.Cells.Item(2, 3).Value = 0.362743313229904     # C2 Estimate
$ws2.Cells.Item(2, 4).Value = 0.331492816698077     # D2 Lower
$ws2.Cells.Item(2, 5).Value = 0.393993809761731     # E2 Upper

$ws2.Cells.Item(3, 3).Value = 0.327066236287387     # C3 Estimate
$ws2.Cells.Item(3, 4).Value = 0.297194139215492     # D3 Lower
$ws2.Cells.Item(3, 5).Value = 0.356938333359282     # E3 Upper

$ws2.Cells.Item(4, 3).Value = 0.411420401256562     # C4 Estimate
$ws2.Cells.Item(4, 4).Value = 0.38117090248048      # D4 Lower
$ws2.Cells.Item(4, 5).Value = 0.441669900032643     # E4 Upper
